# Update computed price/profit columns (H:N) across multiple sheets
# as produced by the scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1661.6666
$ws.Range("I2").Value = 746.7143
$ws.Range("J2").Value = 2462.25
$ws.Range("K2").Value = 746.7143
$ws.Range("L2").Value = 2462.25
$ws.Range("M2").Value = -633.7143
$ws.Range("N2").Value = -2688.25

$ws.Range("H76").Value = 7621.25
$ws.Range("I76").Value = 5990
$ws.Range("J76").Value = 8165
$ws.Range("K76").Value = 5990
$ws.Range("L76").Value = 8165
$ws.Range("M76").Value = -5675
$ws.Range("N76").Value = -8795

$ws.Range("H79").Value = 7621.25
$ws.Range("I79").Value = 5990
$ws.Range("J79").Value = 8165
$ws.Range("K79").Value = 5990
$ws.Range("L79").Value = 8165
$ws.Range("M79").Value = -4898
$ws.Range("N79").Value = -10349

$ws.Range("H80").Value = 1038
$ws.Range("I80").Value = 325
$ws.Range("K80").Value = 975
$ws.Range("M80").Value = 23

$ws.Range("H82").Value = 2639.6
$ws.Range("I82").Value = 1924.5
$ws.Range("K82").Value = 5773.5
$ws.Range("M82").Value = -5367.5

$ws.Range("H83").Value = 1038
$ws.Range("I83").Value = 325
$ws.Range("K83").Value = 2925
$ws.Range("M83").Value = 2067

$ws.Range("H85").Value = 2639.6
$ws.Range("I85").Value = 1924.5
$ws.Range("K85").Value = 5773.5
$ws.Range("M85").Value = -4369.5

$ws.Range("H103").Value = 999.75
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws.Range("H138").Value = 3358.3333
$ws.Range("I138").Value = 1500
$ws.Range("K138").Value = 4500
$ws.Range("M138").Value = 640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4442.8335
$ws.Range("I110").Value = 2771.6667
$ws.Range("K110").Value = 2771.6667
$ws.Range("M110").Value = -726.6667000000002

$ws.Range("H122").Value = 2914.6316
$ws.Range("I122").Value = 2909.889
$ws.Range("K122").Value = 8729.667000000001
$ws.Range("M122").Value = -6279.667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2443.25
$ws.Range("I94").Value = 2439.4666
$ws.Range("K94").Value = 2439.4666
$ws.Range("M94").Value = -1988.4666

$ws.Range("H100").Value = 41605
$ws.Range("J100").Value = 41605
$ws.Range("L100").Value = 41605
$ws.Range("N100").Value = -43769

$ws.Range("H134").Value = 7756.6
$ws.Range("I134").Value = 9023.857
$ws.Range("J134").Value = 4799.6665
$ws.Range("K134").Value = 27071.571
$ws.Range("L134").Value = 14398.9995
$ws.Range("M134").Value = -24536.571
$ws.Range("N134").Value = -19468.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3674.1333
$ws.Range("J31").Value = 3374
$ws.Range("L31").Value = 3374
$ws.Range("N31").Value = -3964

$ws.Range("H34").Value = 3674.1333
$ws.Range("J34").Value = 3374
$ws.Range("L34").Value = 3374
$ws.Range("N34").Value = -3778

$ws.Range("H58").Value = 3183.1667
$ws.Range("I58").Value = 2033
$ws.Range("K58").Value = 2033
$ws.Range("M58").Value = -1830

$ws.Range("H136").Value = 3183.1667
$ws.Range("I136").Value = 2033
$ws.Range("K136").Value = 6099
$ws.Range("M136").Value = -3549

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 172.70589
$ws.Range("I2").Value = 238.9
$ws.Range("J2").Value = 78.14286
$ws.Range("K2").Value = 1433.4
$ws.Range("L2").Value = 468.85716
$ws.Range("M2").Value = -1320.4
$ws.Range("N2").Value = -694.85716

$ws.Range("H5").Value = 680.86664
$ws.Range("I5").Value = 515.2143
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 1545.6429
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = -1433.6429
$ws.Range("N5").Value = -9224

$ws.Range("H98").Value = 289
$ws.Range("J98").Value = 289
$ws.Range("L98").Value = 867
$ws.Range("N98").Value = -3863

$ws.Range("H135").Value = 680.86664
$ws.Range("I135").Value = 515.2143
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 4636.928699999999
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -2101.928699999999
$ws.Range("N135").Value = -32070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18836.75
$ws.Range("J11").Value = 24982.334
$ws.Range("L11").Value = 24982.334
$ws.Range("N11").Value = -25260.334

$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("N33").Value = -25504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H22").Value = 1208.1578
$ws.Range("I22").Value = 1035
$ws.Range("J22").Value = 1583.3334
$ws.Range("K22").Value = 1035
$ws.Range("L22").Value = 1583.3334
$ws.Range("M22").Value = -740
$ws.Range("N22").Value = -2173.3334

$ws.Range("H27").Value = 1208.1578
$ws.Range("I27").Value = 1035
$ws.Range("J27").Value = 1583.3334
$ws.Range("K27").Value = 1035
$ws.Range("L27").Value = 1583.3334
$ws.Range("M27").Value = -928
$ws.Range("N27").Value = -1797.3334

$ws.Range("H82").Value = 1423.2307
$ws.Range("I82").Value = 1423.2307
$ws.Range("K82").Value = 1423.2307
$ws.Range("M82").Value = -1062.2307

$ws.Range("H85").Value = 1423.2307
$ws.Range("I85").Value = 1423.2307
$ws.Range("K85").Value = 1423.2307
$ws.Range("M85").Value = -175.2307000000001

$ws.Range("H132").Value = 21127
$ws.Range("I132").Value = 21607.555
$ws.Range("J132").Value = 19397
$ws.Range("K132").Value = 64822.665
$ws.Range("L132").Value = 58191
$ws.Range("M132").Value = -62292.665
$ws.Range("N132").Value = -63251

$ws.Range("H136").Value = 3563.625
$ws.Range("J136").Value = 3495
$ws.Range("L136").Value = 10485
$ws.Range("N136").Value = -15585

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 337.5
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 575
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 575
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -801

$ws.Range("H132").Value = 300
$ws.Range("I132").Value = 300
$ws.Range("K132").Value = 900
$ws.Range("M132").Value = 1630
